# Re-create the view/formatting changes captured in the commit
# "changes for MVC structure of project" for sample_roster.xlsx.
#
# The author simply re-opened the workbook on a new machine/project layout
# (hence the changed absolute path / revision ids - those are environment
# metadata that Excel itself stamps on save and aren't reachable through the
# object model) and, while poking around the sheet, ended up with:
#   - a different zoom level on Sheet1,
#   - a different active cell / selection,
#   - column A widened (no longer auto "best fit").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Zoom Sheet1 to 145% (drives <sheetView zoomScale="145" .../>).
$excel.ActiveWindow.Zoom = 145

# Move the selection/active cell to D2 (drives <selection activeCell="D2" sqref="D2"/>).
$ws.Range("D2").Select()

# Widen column A from the old "best fit" width to a fixed, wider width.
$ws.Columns("A").ColumnWidth = 39.5
